$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data values (price + 1h volume change)
# Columns B/C (name/link) can contain plain text; column D (price) values
# that look like plain numbers must be forced to text so Excel keeps the
# original formatted digits/precision instead of coercing to a float.

$ws.Range('D2').Value = '64.155.02'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '3.620.90'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '403.02'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').Value = '3.771.18'
$ws.Range('E7').Value = '  +6.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.616'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -4.67%  '
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.717'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -7.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -11.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000300'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -8.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.55'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.273.64'
$ws.Range('E14').Value = '  +4.11%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.74'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.138'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.33%  '
$ws.Range('D17').Value = '3.633.05'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.71'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.80'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.07'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.63%  '
$ws.Range('D21').Value = '65.063.71'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '415.31'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -8.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.92'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +15.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.01'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -6.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '35.39'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.15'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.35'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.01'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.26'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.68'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.117'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.54%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.159'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.83'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.81'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.42'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0459'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.92'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +26.86%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.86'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +30.48%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.137'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.27'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.93%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '141.48'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.72%  '
$ws.Range('B45').Value = 'PEPE'
$ws.Range('C45').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D45').Value = '0.0₃0616'
$ws.Range('E45').Value = '  -19.95%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.05'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +17.45%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.02'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.21'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.75'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.51'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.286'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.18%  '
